$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Window position (best effort - moves the window on screen)
$win = $excel.ActiveWindow
$win.Left = 8480
$win.Top = 400

# 2) Extend the B1 comment with the new "4: captain" npcType value
$cmt = $ws.Range("B1").Comment
$cmt.Text($cmt.Text() + "`n4: 船长")

# 3) Add the new shared strings used by the ship-modification dialog rows,
#    in the same order they first appear in the target workbook.
$ws.Range("E10").Value = "name_shipyard_owner"
$ws.Range("A10").Value = "dialog_modify_ship_confirm"
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 0
$ws.Range("F10").Value = 0

$ws.Range("A11").Value = "dialog_no_enough_money"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0

# 4) Move the active selection to E5
$ws.Range("E5").Select()
